{"js": "// Word, Inc. -> the single paragraph's font size is bumped from 16pt\n// (32 half-points) to 36pt (72 half-points), and a trailing run\n// containing a single space (with the same resulting formatting) is\n// appended to the paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n\n// 1) Bump the font size for the paragraph mark + existing run from\n//    16pt to 36pt. `font.size` maps to <w:sz>, `font.sizeBidirectional`\n//    maps to <w:szCs> - both need updating to mirror the diff.\nparagraph.font.size = 36;\nparagraph.font.sizeBidirectional = 36;\nawait context.sync();\n\n// 2) Append a new trailing run holding a single space, with formatting\n//    identical to the (now-resized) existing run, as its own <w:r>.\nconst endRange = paragraph.getRange(\"End\");\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p><w:r><w:rPr>\" +\n  '<w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n  '<w:color w:val=\"4D5156\"/>' +\n  '<w:sz w:val=\"72\"/>' +\n  '<w:szCs w:val=\"72\"/>' +\n  '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>' +\n  '</w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r></w:p>' +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nendRange.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# The document has a single paragraph. The edit:\n#   1) bumps the font size (w:sz / w:szCs) from 16pt (32 half-points) to\n#      36pt (72 half-points) for both the paragraph mark and the existing\n#      run of text, and\n#   2) appends a new trailing run containing a single space \" \", carrying\n#      the same (now-resized) run formatting, as its own sibling <w:r>.\n\n$d = $word.ActiveDocument\n\n$p = $d.Paragraphs(1)\n\n# --- Step 1: resize the font (paragraph mark + run) ----------------------\n# Font.Size writes <w:sz>, Font.SizeBi writes the complex-script\n# counterpart <w:szCs>; both need to be set to mirror the diff.\n$p.Range.Font.Size = 36\n$p.Range.Font.SizeBi = 36\n\n# --- Step 2: add a sibling run holding just a space -----------------------\n# Word's COM object model has no direct \"add a new run\" primitive, so we\n# get there the same way a user typing at the end of the paragraph and\n# pressing Enter then Backspace would: split the paragraph into two (which\n# mints a brand new run, cloned from the original formatting, for the\n# second half), type the space into the new paragraph, then delete the\n# paragraph mark that separates them. Deleting that mark re-merges the two\n# paragraphs back into a single <w:p> while leaving the two runs as\n# distinct sibling elements.\n$endOfText = $p.Range.End\n\n$splitPoint = $d.Range($endOfText, $endOfText)\n$splitPoint.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs($p.Index + 1)\n$newPara.Range.InsertAfter(\" \")\n\n# The paragraph mark that now separates the original text from the new\n# \" \" paragraph sits in the single character just before $endOfText.\n$mark = $d.Range($endOfText - 1, $endOfText)\n$mark.Delete()\n"}
